$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing data (and formatting) to the right.
$ws.Columns("D:D").Insert()

# The freshly inserted column D picks up formatting from column C (to its left).
# Copy the number formatting from column E (the old D column, now shifted one to
# the right) back onto the new column D so the date / number styles match.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 10600
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = 53000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 65200
$ws.Range("D18").Value = -54500
$ws.Range("D20").Value = 3100
$ws.Range("D21").Value = -51400
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -51400
$ws.Range("D24").Value = "NA"
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -51400
$ws.Range("D27").Value = -51400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3100
$ws.Range("D33").Value = -51400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -51400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 49500
$ws.Range("D42").Value = 86700
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 3300
$ws.Range("D46").Value = 139500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 100
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 139900
$ws.Range("D57").Value = 3300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 6300
$ws.Range("D60").Value = 9500
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 9600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -110600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 130400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -51400
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -52700
$ws.Range("D91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -33900
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 125900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 39300
